# Update the month headers: insert a new "Jan_2026" column at D1, shifting
# the existing Dec_2025 / Nov_2025 values one column to the right, and
# dropping the oldest month (Oct_2025) that falls off the end.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Jan_2026"
$ws.Range("E1").Value = "Dec_2025"
$ws.Range("F1").Value = "Nov_2025"
